$d = $word.ActiveDocument

# Locate the sentence that needs the clarifying parenthetical inserted.
$rng = $d.Content
$found = $rng.Find.Execute("Controller is responsible for monitoring the list of ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $startPos = $rng.Start

    $part1 = "Controller is responsible for monitoring "
    $part2 = "(not maintaining. Maintaining the active list is what we discussed in previous lecture) "
    $part3 = "the list of "

    # Rewrite the run's text content in one shot (old run text -> new combined text).
    $rng.Text = $part1 + $part2 + $part3

    $b0 = $startPos
    $b1 = $b0 + $part1.Length
    $b2 = $b1 + $part2.Length
    $b3 = $b2 + $part3.Length

    # Force the engine to keep this as three distinct runs (matching the
    # source document's run layout) instead of silently re-flattening the
    # whole paragraph back into a single run. A harmless formatting
    # round-trip (bold on, bold off) on the two outer pieces pins the run
    # boundaries on either side of the new middle run without altering any
    # visible text formatting (the middle run is left completely alone).
    $left = $d.Range($b0, $b1)
    $left.Bold = 1
    $left.Bold = 0

    $right = $d.Range($b2, $b3)
    $right.Bold = 1
    $right.Bold = 0
}
